# Apply updates to "Project tasks" worksheet per the commit "update minute and task list"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary / minute block (rows 4-7) ---
$ws.Cells.Item(4, 2).Value = 44672.7083333333      # Project finish date
$ws.Cells.Item(5, 2).Value = "19 days"             # Duration
$ws.Cells.Item(6, 2).Value = 0.67                  # % complete
$ws.Cells.Item(7, 2).Value = 44670.4002092593      # Exported on

# --- Task 11 "coding" (row 20): Dependents (after) gains task 13 ---
$ws.Cells.Item(20, 11).Value = "12 - unit testing, 13 - hardware team review"

# --- Task 13 "hardware team review" (row 22): Depends on gains task 11 ---
$ws.Cells.Item(22, 8).Value = "10 - assignment of ports, 11 - coding"

# --- Task 21 "hardware team review" (row 30): Duration/Start changed ---
$ws.Cells.Item(30, 5).Value = "1.67 hours"
$ws.Cells.Item(30, 6).Value = 44655.375

# --- Task 23 "Stepped Motor" (row 32): Duration/Finish changed ---
$ws.Cells.Item(32, 5).Value = "3.94 days"
$ws.Cells.Item(32, 7).Value = 44658.4444444444

# --- Task 30 "code review" (row 39): Start/Finish changed ---
$ws.Cells.Item(39, 6).Value = 44657.6111111111
$ws.Cells.Item(39, 7).Value = 44658.4444444444

# --- Task 47 "Integration of all components" (row 56): Dependents (after) set ---
$ws.Cells.Item(56, 11).Value = "49 - Car Mount"

# --- Task 49 "Car Mount" (row 58): Start/Finish/Depends on changed ---
$ws.Cells.Item(58, 6).Value = 44658.375
$ws.Cells.Item(58, 7).Value = 44659.7083333333
$ws.Cells.Item(58, 8).Value = "48 - Integration of all the interfaces, 47 - Integration of all components"

# --- Task 50 "Report + log book" (row 59): Duration/Start/Finish changed, gains dependent task 51 ---
$ws.Cells.Item(59, 5).Value = "6 days"
$ws.Cells.Item(59, 6).Value = 44662.375
$ws.Cells.Item(59, 7).Value = 44669.7083333333
$ws.Cells.Item(59, 11).Value = "51 - Poster"

# --- New task 51 "Poster" (row 60) ---
$ws.Cells.Item(60, 1).Value = 51
# Outline number column stores numeric-looking values as text (matches column B elsewhere),
# so force text storage with a leading apostrophe.
$ws.Cells.Item(60, 2).Value = "'16"
$ws.Cells.Item(60, 3).Value = "Poster"
$ws.Cells.Item(60, 4).Value = "Cher Khor"
$ws.Cells.Item(60, 5).Value = "3 days"
$ws.Cells.Item(60, 6).Value = 44670.375
$ws.Cells.Item(60, 7).Value = 44672.7083333333
$ws.Cells.Item(60, 8).Value = "50 - Report + log book"
$ws.Cells.Item(60, 9).Value = 0
$ws.Cells.Item(60, 10).Value = "Bucket 1"
$ws.Cells.Item(60, 11).Value = ""
$ws.Cells.Item(60, 12).Value = "24 hours"
$ws.Cells.Item(60, 13).Value = "0 hours"
$ws.Cells.Item(60, 14).Value = "24 hours"
$ws.Cells.Item(60, 15).Value = "No"
$ws.Cells.Item(60, 16).Value = ""
$ws.Cells.Item(60, 17).Value = ""
$ws.Cells.Item(60, 18).Value = ""
$ws.Cells.Item(60, 19).Value = ""

# Copy formatting from the row above (row 59) so number formats / styles match,
# applied after setting values so the text-forced column B keeps style "s=1".
$srcRow = $ws.Range("A59:S59")
$dstRow = $ws.Range("A60:S60")
$srcRow.Copy()
$dstRow.PasteSpecial(-4122)   # xlPasteFormats
